$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 472643.44
$ws.Range("J17").Value = 472643.44
$ws.Range("L17").Value = 1417930.32
$ws.Range("N17").Value = -1418266.32
# Row 103
$ws.Range("H103").Value = 546.10345
$ws.Range("I103").Value = 480.70834
$ws.Range("J103").Value = 860
$ws.Range("K103").Value = 1442.12502
$ws.Range("L103").Value = 2580
$ws.Range("M103").Value = -856.1250199999999
$ws.Range("N103").Value = -3752
# Row 129
$ws.Range("H129").Value = 743.7778
$ws.Range("I129").Value = 649.0769
$ws.Range("K129").Value = 1947.2307
$ws.Range("M129").Value = 3052.7693
# Row 138
$ws.Range("H138").Value = 2096.4937
$ws.Range("I138").Value = 1197.2
$ws.Range("J138").Value = 2401.3389
$ws.Range("K138").Value = 3591.6
$ws.Range("L138").Value = 7204.0167
$ws.Range("M138").Value = 1548.4
$ws.Range("N138").Value = -17484.0167

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1449.25
$ws.Range("I45").Value = 1095
$ws.Range("J45").Value = 1803.5
$ws.Range("K45").Value = 1095
$ws.Range("L45").Value = 1803.5
$ws.Range("M45").Value = -718
$ws.Range("N45").Value = -2557.5
# Row 102
$ws.Range("H102").Value = 2092
$ws.Range("I102").Value = 2092
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2092
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = ""
$ws.Range("N102").Value = -470
# Row 110
$ws.Range("H110").Value = 1820.7858
$ws.Range("I110").Value = 1809.1
$ws.Range("J110").Value = 1850
$ws.Range("K110").Value = 1809.1
$ws.Range("L110").Value = 1850
$ws.Range("M110").Value = 235.9000000000001
$ws.Range("N110").Value = -5940
# Row 122
$ws.Range("H122").Value = 1273.8889
$ws.Range("I122").Value = 1110.8334
$ws.Range("J122").Value = 1600
$ws.Range("K122").Value = 3332.5002
$ws.Range("L122").Value = 4800
$ws.Range("M122").Value = -882.5001999999999
$ws.Range("N122").Value = -9700

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 797164
$ws.Range("I105").Value = 1327723.4
$ws.Range("J105").Value = 1325
$ws.Range("K105").Value = 1327723.4
$ws.Range("L105").Value = 1325
$ws.Range("M105").Value = -1325976.4
$ws.Range("N105").Value = -4819

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 9555.262000000001
$ws.Range("I31").Value = 7486.2812
$ws.Range("J31").Value = 16176
$ws.Range("K31").Value = 7486.2812
$ws.Range("L31").Value = 16176
$ws.Range("M31").Value = -7191.2812
$ws.Range("N31").Value = -16766
# Row 34
$ws.Range("H34").Value = 9555.262000000001
$ws.Range("I34").Value = 7486.2812
$ws.Range("J34").Value = 16176
$ws.Range("K34").Value = 7486.2812
$ws.Range("L34").Value = 16176
$ws.Range("M34").Value = -7284.2812
$ws.Range("N34").Value = -16580
# Row 99
$ws.Range("H99").Value = 3311.7368
$ws.Range("I99").Value = 3343.3076
$ws.Range("J99").Value = 3243.3333
$ws.Range("K99").Value = 3343.3076
$ws.Range("L99").Value = 3243.3333
$ws.Range("M99").Value = -1845.3076
$ws.Range("N99").Value = -6239.3333
# Row 126
$ws.Range("H126").Value = 3311.7368
$ws.Range("I126").Value = 3343.3076
$ws.Range("J126").Value = 3243.3333
$ws.Range("K126").Value = 10029.9228
$ws.Range("L126").Value = 9729.999899999999
$ws.Range("M126").Value = -7559.9228
$ws.Range("N126").Value = -14669.9999
# Row 141
$ws.Range("H141").Value = 32171.28
$ws.Range("J141").Value = 32171.28
$ws.Range("L141").Value = 32171.28
$ws.Range("N141").Value = -42531.28

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 46
$ws.Range("H46").Value = 1984.4348
$ws.Range("I46").Value = 163
$ws.Range("J46").Value = 2257.65
$ws.Range("K46").Value = 489
$ws.Range("L46").Value = 6772.950000000001
$ws.Range("M46").Value = -398
$ws.Range("N46").Value = -6954.950000000001
# Row 59
$ws.Range("H59").Value = 30304362
$ws.Range("J59").Value = 45456044
$ws.Range("L59").Value = 136368132
$ws.Range("N59").Value = -136369212
# Row 60
$ws.Range("H60").Value = 695
$ws.Range("I60").Value = 241
$ws.Range("J60").Value = 1603
$ws.Range("K60").Value = 723
$ws.Range("L60").Value = 4809
$ws.Range("M60").Value = -472
$ws.Range("N60").Value = -5311
# Row 87
$ws.Range("H87").Value = 40522.656
$ws.Range("I87").Value = 2471.3
$ws.Range("K87").Value = 7413.900000000001
$ws.Range("M87").Value = -6165.900000000001
# Row 90
$ws.Range("H90").Value = 40522.656
$ws.Range("I90").Value = 2471.3
$ws.Range("K90").Value = 22241.7
$ws.Range("M90").Value = -16001.7
# Row 113
$ws.Range("H113").Value = 441.17392
$ws.Range("I113").Value = 412.33334
$ws.Range("J113").Value = 472.63635
$ws.Range("K113").Value = 1237.00002
$ws.Range("L113").Value = 1417.90905
$ws.Range("M113").Value = 932.9999800000001
$ws.Range("N113").Value = -5757.90905
# Row 131
$ws.Range("H131").Value = 754.8081
$ws.Range("I131").Value = 496.66666
$ws.Range("J131").Value = 771.46234
$ws.Range("K131").Value = 1489.99998
$ws.Range("L131").Value = 2314.38702
$ws.Range("M131").Value = 3550.00002
$ws.Range("N131").Value = -12394.38702

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 23487.889
$ws.Range("I102").Value = 15501.429
$ws.Range("J102").Value = 28570.182
$ws.Range("K102").Value = 15501.429
$ws.Range("L102").Value = 28570.182
$ws.Range("M102").Value = -13879.429
$ws.Range("N102").Value = -31814.182
# Row 122
$ws.Range("H122").Value = 1270.7142
$ws.Range("I122").Value = 1282.5
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 3847.5
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -1397.5
$ws.Range("N122").Value = -8500

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2554.2222
$ws.Range("I7").Value = 1998.5883
$ws.Range("J7").Value = 12000
$ws.Range("K7").Value = 1998.5883
$ws.Range("L7").Value = 12000
$ws.Range("M7").Value = -1886.5883
$ws.Range("N7").Value = -12224
# Row 40
$ws.Range("H40").Value = 59900.117
$ws.Range("I40").Value = 1058.1666
$ws.Range("J40").Value = 201120.8
$ws.Range("K40").Value = 1058.1666
$ws.Range("L40").Value = 201120.8
$ws.Range("M40").Value = -922.1666
$ws.Range("N40").Value = -201392.8
# Row 82
$ws.Range("H82").Value = 1597.129
$ws.Range("I82").Value = 1260.421
$ws.Range("J82").Value = 2130.25
$ws.Range("K82").Value = 1260.421
$ws.Range("L82").Value = 2130.25
$ws.Range("M82").Value = -899.421
$ws.Range("N82").Value = -2852.25
# Row 85
$ws.Range("H85").Value = 1597.129
$ws.Range("I85").Value = 1260.421
$ws.Range("J85").Value = 2130.25
$ws.Range("K85").Value = 1260.421
$ws.Range("L85").Value = 2130.25
$ws.Range("M85").Value = -12.42100000000005
$ws.Range("N85").Value = -4626.25
# Row 122
$ws.Range("H122").Value = 1600
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1600
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = ""
$ws.Range("M122").Value = 4800
$ws.Range("N122").Value = -9700
# Row 126
$ws.Range("H126").Value = 2554.2222
$ws.Range("I126").Value = 1998.5883
$ws.Range("J126").Value = 12000
$ws.Range("K126").Value = 5995.7649
$ws.Range("L126").Value = 36000
$ws.Range("M126").Value = -3525.7649
$ws.Range("N126").Value = -40940

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1153.9
$ws.Range("I96").Value = 904.58826
$ws.Range("J96").Value = 2566.6667
$ws.Range("K96").Value = 904.58826
$ws.Range("L96").Value = 2566.6667
$ws.Range("M96").Value = 468.41174
$ws.Range("N96").Value = -5312.6667
# Row 122
$ws.Range("H122").Value = 8243.467000000001
$ws.Range("I122").Value = 3825.6667
$ws.Range("J122").Value = 11188.667
$ws.Range("K122").Value = 11477.0001
$ws.Range("L122").Value = 33566.001
$ws.Range("M122").Value = -9027.000100000001
$ws.Range("N122").Value = -38466.001
# Row 126
$ws.Range("H126").Value = 863.86664
$ws.Range("I126").Value = 875.2857
$ws.Range("J126").Value = 704
$ws.Range("K126").Value = 2625.8571
$ws.Range("L126").Value = 2112
$ws.Range("M126").Value = -155.8571000000002
$ws.Range("N126").Value = -7052
